$wb = $excel.ActiveWorkbook

$wsWeo  = $wb.Worksheets.Item("weo_pg")
$wsCcs  = $wb.Worksheets.Item("ccs_retrofits")
$wsStk  = $wb.Worksheets.Item("existing_stock")

# ---------------------------------------------------------------------------
# 1) Title banner (row 1) font size 8 -> 9 pt on all three sheets
# ---------------------------------------------------------------------------
$wsWeo.Range("A1:Q1").Font.Size = 9
$wsCcs.Range("A1:U1").Font.Size = 9
$wsStk.Range("A1:T1").Font.Size = 9

# ---------------------------------------------------------------------------
# 2) Lower-case the two commodity labels on the existing_stock sheet
# ---------------------------------------------------------------------------
$wsStk.Range("C66:C81").Value = "elc_sol-BGR"
$wsStk.Range("C82:C88").Value = "elc_win-BGR"

# ---------------------------------------------------------------------------
# 3) ccs_retrofits sheet: rewrite rows 27-33 and append new rows 34-48
# ---------------------------------------------------------------------------

# -- rows 27-29 become the "ep_coal_subcritical__m_ccs-rf / coal" record --
$wsCcs.Range("A27").Value = "ep_coal_subcritical__m_ccs-rf"
$wsCcs.Range("D27").Value = [double]"9.3936000000000006E-2"
$wsCcs.Range("I27").Value = "ep_coal_subcritical__m"

$wsCcs.Range("A28").Value = "ep_coal_subcritical__m_ccs-rf"
$wsCcs.Range("B28").Value = "coal"
$wsCcs.Range("D28").Value = [double]"9.3936000000000006E-2"
$wsCcs.Range("E28").Value = 3583
$wsCcs.Range("F28").Value = 49
$wsCcs.Range("G28").Value = [double]"4.59"
$wsCcs.Range("H28").Value = [double]"0.5605"
$wsCcs.Range("I28").Value = "ep_coal_subcritical__m"

$wsCcs.Range("A29").Value = "ep_coal_subcritical__m_ccs-rf"
$wsCcs.Range("B29").Value = "coal"
$wsCcs.Range("D29").Value = [double]"9.3936000000000006E-2"
$wsCcs.Range("E29").Value = 3583
$wsCcs.Range("F29").Value = 49
$wsCcs.Range("G29").Value = [double]"4.59"
$wsCcs.Range("H29").Value = [double]"0.5605"
$wsCcs.Range("I29").Value = "ep_coal_subcritical__m"

# -- rows 30-33 become the "ep_coal_subcritical_ccs-rf / coal" record --
$wsCcs.Range("A30").Value = "ep_coal_subcritical_ccs-rf"
$wsCcs.Range("B30").Value = "coal"
$wsCcs.Range("D30").Value = [double]"0.10646080000000004"
$wsCcs.Range("E30").Value = 3583
$wsCcs.Range("F30").Value = 49
$wsCcs.Range("G30").Value = [double]"4.59"
$wsCcs.Range("H30").Value = [double]"0.5605"
$wsCcs.Range("I30").Value = "ep_coal_subcritical"

$wsCcs.Range("A31").Value = "ep_coal_subcritical_ccs-rf"
$wsCcs.Range("B31").Value = "coal"
$wsCcs.Range("D31").Value = [double]"0.10646080000000004"
$wsCcs.Range("E31").Value = 3583
$wsCcs.Range("F31").Value = 49
$wsCcs.Range("G31").Value = [double]"4.59"
$wsCcs.Range("H31").Value = [double]"0.5605"
$wsCcs.Range("I31").Value = "ep_coal_subcritical"

$wsCcs.Range("A32").Value = "ep_coal_subcritical_ccs-rf"
$wsCcs.Range("B32").Value = "coal"
$wsCcs.Range("D32").Value = [double]"0.10646080000000004"
$wsCcs.Range("E32").Value = 3583
$wsCcs.Range("F32").Value = 49
$wsCcs.Range("G32").Value = [double]"4.59"
$wsCcs.Range("H32").Value = [double]"0.5605"
$wsCcs.Range("I32").Value = "ep_coal_subcritical"

$wsCcs.Range("A33").Value = "ep_coal_subcritical_ccs-rf"
$wsCcs.Range("B33").Value = "coal"
$wsCcs.Range("D33").Value = [double]"0.10646080000000004"
$wsCcs.Range("E33").Value = 3583
$wsCcs.Range("F33").Value = 49
$wsCcs.Range("G33").Value = [double]"4.59"
$wsCcs.Range("H33").Value = [double]"0.5605"
$wsCcs.Range("I33").Value = "ep_coal_subcritical"

# -- stamp the banded (white/grey) row formatting down onto the new rows --
$wsCcs.Range("A32:K33").Copy() | Out-Null
$wsCcs.Range("A34:K48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# -- rows 34-36: ep_coal_subcritical_ccs-rf / coal (3 more copies) --
$coalRows = 34,35,36
foreach ($r in $coalRows) {
    $wsCcs.Range("A$r").Value = "ep_coal_subcritical_ccs-rf"
    $wsCcs.Range("B$r").Value = "coal"
    $wsCcs.Range("C$r").Value = "ELC"
    $wsCcs.Range("D$r").Value = [double]"0.10646080000000004"
    $wsCcs.Range("E$r").Value = 3583
    $wsCcs.Range("F$r").Value = 49
    $wsCcs.Range("G$r").Value = [double]"4.59"
    $wsCcs.Range("H$r").Value = [double]"0.5605"
    $wsCcs.Range("I$r").Value = "ep_coal_subcritical"
    $wsCcs.Range("J$r").Value = 1
    $wsCcs.Range("K$r").Value = 20
}

# -- rows 37-39: ep_gas_combined_cycle_ccs-rf / gas --
$gasRows = 37,38,39
foreach ($r in $gasRows) {
    $wsCcs.Range("A$r").Value = "ep_gas_combined_cycle_ccs-rf"
    $wsCcs.Range("B$r").Value = "gas"
    $wsCcs.Range("C$r").Value = "ELC"
    $wsCcs.Range("D$r").Value = [double]"0.50470000000000004"
    $wsCcs.Range("E$r").Value = 1365
    $wsCcs.Range("F$r").Value = [double]"34.200000000000003"
    $wsCcs.Range("G$r").Value = [double]"3.75"
    $wsCcs.Range("H$r").Value = [double]"0.84455000000000002"
    $wsCcs.Range("I$r").Value = "ep_gas_combined_cycle"
    $wsCcs.Range("J$r").Value = 1
    $wsCcs.Range("K$r").Value = 20
}

# -- row 40: ep_gas_steam_turbine_G100000400048_ccs-rf / gas --
$wsCcs.Range("A40").Value = "ep_gas_steam_turbine_G100000400048_ccs-rf"
$wsCcs.Range("B40").Value = "gas"
$wsCcs.Range("C40").Value = "ELC"
$wsCcs.Range("D40").Value = [double]"0.27037500000000003"
$wsCcs.Range("E40").Value = 1365
$wsCcs.Range("F40").Value = [double]"34.200000000000003"
$wsCcs.Range("G40").Value = [double]"3.75"
$wsCcs.Range("H40").Value = [double]"0.84455000000000002"
$wsCcs.Range("I40").Value = "ep_gas_steam_turbine_G100000400048"
$wsCcs.Range("J40").Value = 1
$wsCcs.Range("K40").Value = 20

# -- row 41: ep_gas_steam_turbine_G100000406325_ccs-rf / gas --
$wsCcs.Range("A41").Value = "ep_gas_steam_turbine_G100000406325_ccs-rf"
$wsCcs.Range("B41").Value = "gas"
$wsCcs.Range("C41").Value = "ELC"
$wsCcs.Range("D41").Value = [double]"0.31093124999999994"
$wsCcs.Range("E41").Value = 1365
$wsCcs.Range("F41").Value = [double]"34.200000000000003"
$wsCcs.Range("G41").Value = [double]"3.75"
$wsCcs.Range("H41").Value = [double]"0.84455000000000002"
$wsCcs.Range("I41").Value = "ep_gas_steam_turbine_G100000406325"
$wsCcs.Range("J41").Value = 1
$wsCcs.Range("K41").Value = 20

# -- row 42: ep_gas_steam_turbine_G100000406326__m_ccs-rf / gas --
$wsCcs.Range("A42").Value = "ep_gas_steam_turbine_G100000406326__m_ccs-rf"
$wsCcs.Range("B42").Value = "gas"
$wsCcs.Range("C42").Value = "ELC"
$wsCcs.Range("D42").Value = [double]"0.31093124999999994"
$wsCcs.Range("E42").Value = 1365
$wsCcs.Range("F42").Value = [double]"34.200000000000003"
$wsCcs.Range("G42").Value = [double]"3.75"
$wsCcs.Range("H42").Value = [double]"0.84455000000000002"
$wsCcs.Range("I42").Value = "ep_gas_steam_turbine_G100000406326__m"
$wsCcs.Range("J42").Value = 1
$wsCcs.Range("K42").Value = 20

# -- row 43: ep_gas_steam_turbine_G100000406327__m_ccs-rf / gas --
$wsCcs.Range("A43").Value = "ep_gas_steam_turbine_G100000406327__m_ccs-rf"
$wsCcs.Range("B43").Value = "gas"
$wsCcs.Range("C43").Value = "ELC"
$wsCcs.Range("D43").Value = [double]"0.31093124999999994"
$wsCcs.Range("E43").Value = 1365
$wsCcs.Range("F43").Value = [double]"34.200000000000003"
$wsCcs.Range("G43").Value = [double]"3.75"
$wsCcs.Range("H43").Value = [double]"0.84455000000000002"
$wsCcs.Range("I43").Value = "ep_gas_steam_turbine_G100000406327__m"
$wsCcs.Range("J43").Value = 1
$wsCcs.Range("K43").Value = 20

# -- rows 44-48: ep_gas_steam_turbine_ccs-rf / gas (5 copies) --
$turbineRows = 44,45,46,47,48
foreach ($r in $turbineRows) {
    $wsCcs.Range("A$r").Value = "ep_gas_steam_turbine_ccs-rf"
    $wsCcs.Range("B$r").Value = "gas"
    $wsCcs.Range("C$r").Value = "ELC"
    $wsCcs.Range("D$r").Value = [double]"0.27938750000000001"
    $wsCcs.Range("E$r").Value = 1365
    $wsCcs.Range("F$r").Value = [double]"34.200000000000003"
    $wsCcs.Range("G$r").Value = [double]"3.75"
    $wsCcs.Range("H$r").Value = [double]"0.84455000000000002"
    $wsCcs.Range("I$r").Value = "ep_gas_steam_turbine"
    $wsCcs.Range("J$r").Value = 1
    $wsCcs.Range("K$r").Value = 20
}
